$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (cell styles) of the last existing row (1230) down into
# the 33 new rows (1231:1263) before writing values, so the new cells pick up
# the same styles (date style on A, number style on B/C, percent style on D)
# without creating any new style entries.
$ws.Range("A1230:D1230").Copy()
$ws.Range("A1231:D1263").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Build the 33x4 block of new values (columns A:D, rows 1231:1263) and write
# it in one shot.
$data = New-Object 'object[,]' 33,4
$data[0,0] = 45908
$data[0,1] = 164.10166062618396
$data[0,2] = 143.14565198034285
$data[0,3] = 0.13189999999999999
$data[1,0] = 45909
$data[1,1] = 165.08807713855271
$data[1,2] = 142.97070710680612
$data[1,3] = 0.13214999999999999
$data[2,0] = 45910
$data[2,1] = 166.19277902183612
$data[2,2] = 143.70809232857292
$data[2,3] = 0.13200000000000001
$data[3,0] = 45911
$data[3,1] = 169.18682260202965
$data[3,2] = 144.51789255281409
$data[3,3] = 0.13214999999999999
$data[4,0] = 45912
$data[4,1] = 167.52217638921519
$data[4,2] = 143.63023585302813
$data[4,3] = 0.1321
$data[5,0] = 45915
$data[5,1] = 167.9499629109979
$data[5,2] = 144.9174117648484
$data[5,3] = 0.13159999999999999
$data[6,0] = 45916
$data[6,1] = 168.85643578386353
$data[6,2] = 145.43749140620787
$data[6,3] = 0.13070000000000001
$data[7,0] = 45917
$data[7,1] = 170.40905000402518
$data[7,2] = 146.98401054939092
$data[7,3] = 0.1305
$data[8,0] = 45918
$data[8,1] = 170.10457517279448
$data[8,2] = 146.8889715373605
$data[8,3] = 0.13105
$data[9,0] = 45919
$data[9,1] = 171.40267484404546
$data[9,2] = 147.25808311138383
$data[9,3] = 0.13144999999999998
$data[10,0] = 45922
$data[10,1] = 167.02870885913737
$data[10,2] = 146.49500484886738
$data[10,3] = 0.1326
$data[11,0] = 45923
$data[11,1] = 170.70365189075872
$data[11,2] = 147.82325933939509
$data[11,3] = 0.13140000000000002
$data[12,0] = 45924
$data[12,1] = 171.44937494968545
$data[12,2] = 147.89070735717448
$data[12,3] = 0.13175000000000001
$data[13,0] = 45925
$data[13,1] = 167.57217774462723
$data[13,2] = 146.69386595562062
$data[13,3] = 0.13250000000000001
$data[14,0] = 45926
$data[14,1] = 170.23544944641685
$data[14,2] = 146.83563702487311
$data[14,3] = 0.13225000000000001
$data[15,0] = 45929
$data[15,1] = 169.545349622374
$data[15,2] = 147.73427762577322
$data[15,3] = 0.13289999999999999
$data[16,0] = 45930
$data[16,1] = 166.7104591396666
$data[16,2] = 147.63354475323877
$data[16,3] = 0.13255
$data[17,0] = 45931
$data[17,1] = 165.66983641408905
$data[17,2] = 146.90700209562334
$data[17,3] = 0.13255
$data[18,0] = 45932
$data[18,1] = 160.31579457539459
$data[18,2] = 145.32432088094117
$data[18,3] = 0.13390000000000002
$data[19,0] = 45933
$data[19,1] = 161.31211483707949
$data[19,2] = 145.57772795986352
$data[19,3] = 0.13425000000000001
$data[20,0] = 45936
$data[20,1] = 159.37654537508769
$data[20,2] = 144.9794990735362
$data[20,3] = 0.13369999999999999
$data[21,0] = 45937
$data[21,1] = 152.06972958297669
$data[21,2] = 142.7063464132616
$data[21,3] = 0.1346
$data[22,0] = 45938
$data[22,1] = 151.28330308984656
$data[22,2] = 143.5028306765013
$data[22,3] = 0.13414999999999999
$data[23,0] = 45939
$data[23,1] = 150.74210349269146
$data[23,2] = 143.06146562796113
$data[23,3] = 0.13339999999999999
$data[24,0] = 45940
$data[24,1] = 149.37039057224675
$data[24,2] = 142.02379993308705
$data[24,3] = 0.1341
$data[25,0] = 45943
$data[25,1] = 149.12824929231064
$data[25,2] = 143.13735348152309
$data[25,3] = 0.1341
$data[26,0] = 45944
$data[26,1] = 147.82037935929264
$data[26,2] = 143.03602497464516
$data[26,3] = 0.1336
$data[27,0] = 45945
$data[27,1] = 150.15682219951336
$data[27,2] = 143.96548712894756
$data[27,3] = 0.13305
$data[28,0] = 45946
$data[28,1] = 147.74757730366653
$data[28,2] = 143.55799247400861
$data[28,3] = 0.1333
$data[29,0] = 45947
$data[29,1] = 147.37483899347922
$data[29,2] = 144.76804888159052
$data[29,3] = 0.1333
$data[30,0] = 45950
$data[30,1] = 149.86941837215835
$data[30,2] = 145.88934567649221
$data[30,3] = 0.13225000000000001
$data[31,0] = 45951
$data[31,1] = 148.98786659788291
$data[31,2] = 145.46112496549864
$data[31,3] = 0.13220000000000001
$data[32,0] = 45952
$data[32,1] = 152.28368445350608
$data[32,2] = 146.25628671858578
$data[32,3] = 0.13144999999999998

$ws.Range("A1231:D1263").Value = $data
